# Adjustments to scenario 3a and 3b parameters for STH
#
# "Platform Coverage" sheet: the old row 3 (Trichuris / MDA / Campaign /
# min-age 0 / F=0,G=15 with a 0.7 coverage series in P..AD) is deleted
# outright, shifting rows 4-8 up to become rows 3-7. Row 2's 0.6 coverage
# series (previously only in H/J/L/N) is then extended with matching 0.6
# values all the way out to column AD (P,R,T,V,X,Z,AB,AD), mirroring the
# pattern that used to live on the deleted row.
#
# "MarketShare" sheet: New Product A's 100% share (row 2) is removed for
# 2026 onward (columns L:Z), and Old Product B / SOC's share (row 3), which
# previously only ran through 2025 (column K), is extended with 100% share
# through 2040 (columns L:Z) to replace it.

$wb = $excel.ActiveWorkbook

# ---- Platform Coverage ----
$ws1 = $wb.Worksheets.Item("Platform Coverage")

# Delete the old row 3 entirely; rows below shift up by one.
$ws1.Rows.Item(3).Delete()

# Extend row 2's coverage series out to 2040 at the same 0.6 level.
$ws1.Range("P2").Value = 0.6
$ws1.Range("R2").Value = 0.6
$ws1.Range("T2").Value = 0.6
$ws1.Range("V2").Value = 0.6
$ws1.Range("X2").Value = 0.6
$ws1.Range("Z2").Value = 0.6
$ws1.Range("AB2").Value = 0.6
$ws1.Range("AD2").Value = 0.6

# ---- MarketShare ----
$ws2 = $wb.Worksheets.Item("MarketShare")

# New Product A loses its share for 2026-2040.
$ws2.Range("L2:Z2").ClearContents()

# Old Product B (SOC) picks up 100% share for 2026-2040.
$ws2.Range("L3:Z3").Value = 1

# ---- Restore on-screen selection/active-sheet state ----
$ws2.Activate()
$ws2.Range("Z3").Select()

$ws1.Activate()
$ws1.Range("G6").Select()
